# Add a new "DefaultHitTime" property row to the Skill.xlsx "Property" sheet.
# (commit message: "add default hit time")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 31 describing the DefaultHitTime property, following the
# same shape as the existing rows (Id / Type / Public / Private / Save /
# View / Index / SaveInterval / RelationValue / Desc).
$ws.Range("A31").Value = "DefaultHitTime"
$ws.Range("B31").Value = "float"
$ws.Range("C31").Value = $false
$ws.Range("D31").Value = $false
$ws.Range("E31").Value = $false
$ws.Range("F31").Value = $true
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = "Friend"
$ws.Range("J31").Value = "缺省打击时间（本来应该打到但是物理没碰撞到或者其他原因）"

# Match the text-number-format styling ("s=1") used by the rest of the table
# for the text columns (A, B, I, J) of this row.
$ws.Range("A31").NumberFormat = "@"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("I31").NumberFormat = "@"
$ws.Range("J31").NumberFormat = "@"

# Leave the selection on the cell past the new row, matching the author's
# cursor position after typing the new row in.
$ws.Range("J32").Select()
